$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.331.18"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "'2.289.02"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'322.86"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "'103.33"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.606"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "'39.87"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").Value = "'0.0908"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "'8.37"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "'0.969"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "'15.19"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "'2.635.92"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "'2.289.39"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'42.313.92"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'7.38"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "'13.43"
$ws.Range("E21").Value = "  +33.23%  "
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("D23").Value = "'73.24"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("D24").Value = "'269.44"
$ws.Range("E24").Value = "  -5.78%  "
$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "'10.91"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("E29").Value = "  -3.74%  "
$ws.Range("D30").Value = "'38.03"
$ws.Range("E30").Value = "  +7.13%  "
$ws.Range("D31").Value = "'164.61"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "'6.15"
$ws.Range("E32").Value = "  +4.02%  "
$ws.Range("D33").Value = "'0.0880"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("E36").Value = "  -14.07%  "
$ws.Range("D37").Value = "'4.62"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("D38").Value = "'0.0355"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("D39").Value = "'3.72"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("D40").Value = "'2.74"
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("D42").Value = "'69.54"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "'93.45"
$ws.Range("E45").Value = "  -8.80%  "
$ws.Range("D46").Value = "'12.32"
$ws.Range("E46").Value = "  +1.77%  "
$ws.Range("D47").Value = "'81.25"
$ws.Range("E47").Value = "  +4.14%  "
$ws.Range("D48").Value = "'112.62"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").Value = "'8.94"
$ws.Range("E49").Value = "  -1.98%  "
$ws.Range("D50").Value = "'5.26"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "'1.601.28"
$ws.Range("E51").Value = "  +3.03%  "
